# Slide 1, "TextBox 4" shape: turn the trailing GitHub URL into a real
# hyperlink, split across two runs (": " stays plain text; the URL is
# split into "https://github.com" + "/miju1234", both pointing at the
# same external hyperlink).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 4")
$tr = $shape.TextFrame.TextRange

$url = "https://github.com/miju1234"

# Full current text is: "깃허브 주소 : https://github.com/miju1234"
#   chars 1-9   -> "깃허브 주소 " (untouched, already separate runs)
#   chars 10-27 -> "https://github.com"  (18 chars) -> becomes hyperlink run 1
#   chars 28-36 -> "/miju1234"           (9 chars)  -> becomes hyperlink run 2
# (the leading ": " run keeps its own text, shortened from
#  ": https://github.com/miju1234" down to just ": ")

$part1 = $tr.Characters(10, 18)
$part1.ActionSettings(1).Hyperlink.Address = $url

$part2 = $tr.Characters(28, 9)
$part2.ActionSettings(1).Hyperlink.Address = $url
